# Event table and activity diagram
# Changed "Government Level Admin" to "Local Governmet Unit Admin"
# (spelling of "Governmet" kept exactly as used by the source edit)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Event table")

$oldText = "Government Level Admin"
$newText = "Local Governmet Unit Admin"

# Data rows of the event table run from row 2 to row 7, columns B..G
for ($row = 2; $row -le 7; $row++) {
    for ($col = 2; $col -le 7; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $text = $cell.Value2
        if ($text -ne $null -and $text -like "*$oldText*") {
            $cell.Value = $text -replace [regex]::Escape($oldText), $newText
        }
    }
}

# Update the saved view: selection moves to B3 and the sheet no longer
# needs to be scrolled down (topLeftCell back to default / A1).
$ws.Range("A1").Select() | Out-Null
$ws.Range("B3").Select() | Out-Null
